$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at B, shifting existing B..H to C..I
$ws.Columns("B").EntireColumn.Insert()

# Populate the new column B with descriptive labels for rows 1-5
$ws.Range("B1").Value = "principle"
$ws.Range("B2").Value = "real interest rate"
$ws.Range("B3").Value = "number of compounding periods per year"
$ws.Range("B4").Value = "number of years"
$ws.Range("B5").Value = "number of payments"

# Size the new column to fit its (longer) descriptive text content
$ws.Columns("B").ColumnWidth = 34.3

# Update the active selection
$null = $ws.Range("E2").Select()
